$d = $word.ActiveDocument

# The paragraph that starts with "Sodelujete v svetovni aktivnosti..." is
# currently built from dozens of small <w:r> runs interleaved with
# <w:proofErr/> spell-check markers (an artifact of pasting text into Word).
# The commit simplifies it down to a single plain run, and also updates the
# named constellation from "Perseus" to "Bik" (Taurus).
#
# We rebuild the whole <w:p> (keeping its original paragraph attributes and
# <w:pPr>) via Range.InsertXML so that every old run and every stray
# <w:proofErr/> marker inside the paragraph is replaced in one shot.

$targetText = "Sodelujete v svetovni aktivnosti opazovanja in beleženja najšibkejših, s prostim očesom  še vidnih zvezd, kot metode za merjenje svetlobnega onesnaževanja na določenem mestu. Z opazovanjem izbranega ozvezdje Bik na nočnem nebu in s primerjavo videnega z zvezdnimi kartami, se lahko ljudje širom sveta podučijo o tem, kako svetila v njihovem kraju prispevajo k svetlobnemu onesnaževanju.  Vaši prispevki v spletno bazo podatkov bodo pomagali dokumentirati nočno nebo, vidno s prostim očesom."

$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.StartsWith("Sodelujete")) {
        $found = $true

        $frag = '<w:p w14:paraId="65CDA169" w14:textId="7FA013F7" w:rsidR="004615A9" w:rsidRPr="00DB0F3B" w:rsidRDefault="00852C5A" w:rsidP="004615A9" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:pStyle w:val="BasicParagraph"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="-72"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Optima-Regular"/><w:sz w:val="20"/></w:rPr></w:pPr><w:r><w:t>' + $targetText + '</w:t></w:r></w:p>'

        $null = $para.Range.InsertXML($frag)
        break
    }
}

if (-not $found) {
    throw "Could not find the target paragraph starting with 'Sodelujete'"
}
